# Adds height (m) and weight (kg) columns (I, J) for every pokemon row,
# fixes the price figures for Wigglytuff (row 51), and marks cell I53
# with an underline font style (matching the new style added by Excel).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("poke_info")

# --- Fix the price values for Wigglytuff (row 51) ---
$ws.Cells.Item(51, 7).Value = 199
$ws.Cells.Item(51, 8).Value = 189

# --- Add Height (column I) and Weight (column J) for each pokemon row ---
$ws.Cells.Item(2, 9).Value = 2.2
$ws.Cells.Item(2, 10).Value = 135.5
$ws.Cells.Item(3, 9).Value = 1.7
$ws.Cells.Item(3, 10).Value = 53.0
$ws.Cells.Item(4, 9).Value = 2.1
$ws.Cells.Item(4, 10).Value = 360.0
$ws.Cells.Item(5, 9).Value = 2.7
$ws.Cells.Item(5, 10).Value = 225.0
$ws.Cells.Item(6, 9).Value = 2.1
$ws.Cells.Item(6, 10).Value = 135.0
$ws.Cells.Item(7, 9).Value = 1.9
$ws.Cells.Item(7, 10).Value = 52.0
$ws.Cells.Item(8, 9).Value = 1.3
$ws.Cells.Item(8, 10).Value = 187.0
$ws.Cells.Item(9, 9).Value = 1.0
$ws.Cells.Item(9, 10).Value = 34.3
$ws.Cells.Item(10, 9).Value = 1.7
$ws.Cells.Item(10, 10).Value = 90.5
$ws.Cells.Item(11, 9).Value = 1.5
$ws.Cells.Item(11, 10).Value = 132.5
$ws.Cells.Item(12, 9).Value = 2.2
$ws.Cells.Item(12, 10).Value = 210.0
$ws.Cells.Item(13, 9).Value = 2.1
$ws.Cells.Item(13, 10).Value = 85.5
$ws.Cells.Item(14, 9).Value = 0.9
$ws.Cells.Item(14, 10).Value = 26.5
$ws.Cells.Item(15, 9).Value = 1.1
$ws.Cells.Item(15, 10).Value = 33.5
$ws.Cells.Item(16, 9).Value = 2.0
$ws.Cells.Item(16, 10).Value = 82.0
$ws.Cells.Item(17, 9).Value = 1.3
$ws.Cells.Item(17, 10).Value = 26.6
$ws.Cells.Item(18, 9).Value = 1.6
$ws.Cells.Item(18, 10).Value = 52.0
$ws.Cells.Item(19, 9).Value = 0.8
$ws.Cells.Item(19, 10).Value = 14.3
$ws.Cells.Item(20, 9).Value = 1.9
$ws.Cells.Item(20, 10).Value = 95.0
$ws.Cells.Item(21, 9).Value = 1.6
$ws.Cells.Item(21, 10).Value = 48.4
$ws.Cells.Item(22, 9).Value = 1.5
$ws.Cells.Item(22, 10).Value = 40.5
$ws.Cells.Item(23, 9).Value = 0.8
$ws.Cells.Item(23, 10).Value = 25.9
$ws.Cells.Item(24, 9).Value = 2.0
$ws.Cells.Item(24, 10).Value = 42.5
$ws.Cells.Item(25, 9).Value = 1.0
$ws.Cells.Item(25, 10).Value = 25.5
$ws.Cells.Item(26, 9).Value = 0.8
$ws.Cells.Item(26, 10).Value = 24.5
$ws.Cells.Item(27, 9).Value = 1.2
$ws.Cells.Item(27, 10).Value = 54.0
$ws.Cells.Item(28, 9).Value = 1.2
$ws.Cells.Item(28, 10).Value = 24.0
$ws.Cells.Item(29, 9).Value = 1.4
$ws.Cells.Item(29, 10).Value = 42.0
$ws.Cells.Item(30, 9).Value = 1.0
$ws.Cells.Item(30, 10).Value = 25.0
$ws.Cells.Item(31, 9).Value = 1.2
$ws.Cells.Item(31, 10).Value = 180.0
$ws.Cells.Item(32, 9).Value = 1.6
$ws.Cells.Item(32, 10).Value = 550.0
$ws.Cells.Item(33, 9).Value = 6.2
$ws.Cells.Item(33, 10).Value = 162.0
$ws.Cells.Item(34, 9).Value = 0.2
$ws.Cells.Item(34, 10).Value = 0.7
$ws.Cells.Item(35, 9).Value = 0.9
$ws.Cells.Item(35, 10).Value = 4.4
$ws.Cells.Item(36, 9).Value = 0.8
$ws.Cells.Item(36, 10).Value = 36.5
$ws.Cells.Item(37, 9).Value = 0.7
$ws.Cells.Item(37, 10).Value = 21.0
$ws.Cells.Item(38, 9).Value = 1.0
$ws.Cells.Item(38, 10).Value = 20.1
$ws.Cells.Item(39, 9).Value = 1.2
$ws.Cells.Item(39, 10).Value = 22.2
$ws.Cells.Item(40, 9).Value = 1.7
$ws.Cells.Item(40, 10).Value = 52.2
$ws.Cells.Item(41, 9).Value = 2.5
$ws.Cells.Item(41, 10).Value = 200.5
$ws.Cells.Item(42, 9).Value = 3.3
$ws.Cells.Item(42, 10).Value = 63.0
$ws.Cells.Item(43, 9).Value = 1.0
$ws.Cells.Item(43, 10).Value = 108.0
$ws.Cells.Item(44, 9).Value = 1.0
$ws.Cells.Item(44, 10).Value = 23.5
$ws.Cells.Item(45, 9).Value = 1.5
$ws.Cells.Item(45, 10).Value = 38.0
$ws.Cells.Item(46, 9).Value = 1.7
$ws.Cells.Item(46, 10).Value = 79.5
$ws.Cells.Item(47, 9).Value = 2.0
$ws.Cells.Item(47, 10).Value = 202.0
$ws.Cells.Item(48, 9).Value = 1.0
$ws.Cells.Item(48, 10).Value = 27.0
$ws.Cells.Item(49, 9).Value = 1.0
$ws.Cells.Item(49, 10).Value = 29.0
$ws.Cells.Item(50, 9).Value = 2.0
$ws.Cells.Item(50, 10).Value = 100.0
$ws.Cells.Item(51, 9).Value = 1.0
$ws.Cells.Item(51, 10).Value = 12.0
$ws.Cells.Item(52, 9).Value = 1.9
$ws.Cells.Item(52, 10).Value = 51.5

# --- New underlined style used on (empty) cell I53 ---
$ws.Range("I53").Font.Underline = 2

# --- Keep the trailing blank rows present in the sheet ---
$ws.Rows(54).OutlineLevel = 0
$ws.Rows(55).OutlineLevel = 0
$ws.Rows(56).OutlineLevel = 0

# --- Update the view: selection & scroll position ---
$ws.Activate()
$ws.Range("I53").Select()
$excel.ActiveWindow.ScrollRow = 50
$excel.ActiveWindow.ScrollColumn = 1
